# Insert a new row at position 400, shifting existing rows 400-495 down to 401-496,
# then populate the new row 400 with the new data record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(400).Insert()

$ws.Cells.Item(400, 1).Value = 10
$ws.Cells.Item(400, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(400, 3).Value = "La Araucanía"
$ws.Cells.Item(400, 4).Value = 44508
$ws.Cells.Item(400, 5).Value = 9
$ws.Cells.Item(400, 6).Value = "Fruta"
$ws.Cells.Item(400, 7).Value = 100104
$ws.Cells.Item(400, 8).Value = "Frutos de pepita"
$ws.Cells.Item(400, 9).Value = 100104005
$ws.Cells.Item(400, 10).Value = "Pera"
$ws.Cells.Item(400, 11).Value = "Packham's Triumph"
$ws.Cells.Item(400, 12).Value = "Primera"
$ws.Cells.Item(400, 13).Value = 95
$ws.Cells.Item(400, 14).Value = 14000
$ws.Cells.Item(400, 15).Value = 14000
$ws.Cells.Item(400, 16).Value = 14000
$ws.Cells.Item(400, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(400, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(400, 19).Value = 778
$ws.Cells.Item(400, 20).Value = 18
